# edit.ps1 - reorders LOB1053.docx sections per the target diff.
# Strategy: paragraphs are addressed structurally (Paragraphs.Item(n)) wherever the
# destination paragraph is uniquely identified by position, and via Find/Replace (with
# collision-safe placeholder tokens) wherever multiple runs inside ONE paragraph must be
# rotated among each other.

$d = $word.ActiveDocument

# --- Paragraph 6 (Objetivos, PT, plain run): APRESENTAR -> CARGA_SHORT ---
$d.Paragraphs.Item(6).Range.Text = "Carga e força elétrica, Campo elétrico, Lei de Gauss, Potencial elétrico, Capacitores e dielétricos, Corrente e Resistência elétrica, Campo magnético: Lei de Biot-Savart, Lei de Ampère Indução eletromagnética e indutância: Lei de Faraday, lei de Lenz, Propriedades magnéticas da matéria, Equações de Maxwell."

# --- Paragraph 7 (Objetivos, EN, italic run): TO_INTRODUCE -> ELECTRIC_SHORT ---
$d.Paragraphs.Item(7).Range.Text = "Electric Charge and Matter. Electric fields. The Gauss' law . Electric Potential . Capacitors and Dielectrics. Electric Current and Resistance. Magnetic Fields . Magnetic Fields sources. Electromagnetic induction and inductance . Magnetic Properties of Matter. Maxwell's equations."

# --- Paragraph 9 (Docente(s) Responsavel(eis), ListBullet run): BERTHA -> APRESENTAR ---
$d.Paragraphs.Item(9).Range.Text = "Apresentar aos estudantes os conceitos básicos do eletromagnetismo tais como carga elétrica, campo elétrico, potencial elétrico, campo magnético e força de Lorentz, mostrando suas aplicações a vários dispositivos e configurações. Adicionalmente, os estudantes irão se familiarizar com as leis de Gauss, Ampère e Faraday. Finalmente, os estudantes devem entender a relação entre campos magnéticos e elétricos e como gerar corrente elétrica apartir de um campo magnético através da indução"

# --- Paragraph 11 (Programa resumido, PT, plain run): CARGA_SHORT -> PROGRAMA_PT_FULL (11 lines) ---
$d.Paragraphs.Item(11).Range.Text = "1) Carga e Força elétrica: carga elétrica; condutores e isolantes; lei de Coulomb; quantização e conservação de cargas." + [char]11 + "2) Campo Elétrico: conceito; linhas de campo; carga pontual e dipolo elétrico, distribuição contínua." + [char]11 + "3) A Lei de Gauss: fluxo; aplicações em simetrias cilíndricas, planares e esféricas." + [char]11 + "4) Potencial Elétrico: conceito e cálculo; energia, potencial e campo elétrico, superfícies equipotenciais; carga puntiforme, dipolo elétrico e distribuições contínuas." + [char]11 + "5) Capacitores e Dielétricos: capacitância, energia e cálculo; associações, dielétrico." + [char]11 + "6) Corrente e Resistência Elétrica: corrente e densidade, resistência, Resistividade e Condutividade em função da temperatura; lei de Ohm, potência, semicondutores e supercondutores." + [char]11 + "7) Campos Magnéticos: lei de Biot-Savart." + [char]11 + "8) Lei de Ampère e aplicações; campo magnético de uma espira, solenoide e toroides." + [char]11 + "9) Indução Eletromagnética: conceitos; Lei de indução de Faraday; Lei de Lenz;" + [char]11 + "10) Propriedades magnéticas da matéria;" + [char]11 + "11) Equações de Maxwell."

# --- Paragraph 12 (Programa resumido, EN, italic run): ELECTRIC_SHORT -> TO_INTRODUCE ---
$d.Paragraphs.Item(12).Range.Text = "To introduce to students the basic concepts of electromagnetism such as electric charge, electric field, electric potential, magnetic field, and Lorentz force showing their applications to several devices and configurations. In addition, the students are going to get familiarized with Gauss, Ampère, and Faraday laws. Finally, students should understand the relation between magnetic and electric fields and how to generate electric current from a magnetic field through induction."

# --- Paragraph 14 (Programa, PT, plain run): PROGRAMA_PT_FULL -> NF_EQUALS (1 line) ---
$d.Paragraphs.Item(14).Range.Text = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# --- Paragraph 17 (Avaliacao: Metodo / Criterio / Norma de recuperacao) ---
# Phase 1: move each current value to a unique placeholder token to avoid collisions
# (the same literal text, e.g. "NF>= 5,0.", needs to land in two different spots).
$rng = $d.Content
$null = $rng.Find.Execute("NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_METODO@@", 2)
$rng = $d.Content
$null = $rng.Find.Execute("NF≥ 5,0.", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_CRITERIO@@", 2)
$rng = $d.Content
$null = $rng.Find.Execute("(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_NORMA@@", 2)

# Phase 2: write the final values into each placeholder location.
$rng = $d.Content
$null = $rng.Find.Execute("@@PH_METODO@@", $true, $false, $false, $false, $false, $true, 1, $false, "NF≥ 5,0.", 2)
$rng = $d.Content
$null = $rng.Find.Execute("@@PH_CRITERIO@@", $true, $false, $false, $false, $false, $true, 1, $false, "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.", 2)
$rng = $d.Content
$null = $rng.Find.Execute("@@PH_NORMA@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 3, Edgard Blucher (2008)." + [char]11 + "RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.3, LTC (2008)." + [char]11 + "TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.3, LTC (2008)." + [char]11 + "SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 3, Pearson Addison Wesley (2009)." + [char]11 + "JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 3, Thomson Pioneira (2008)."

# --- Paragraph 19 (Bibliografia list, plain run, multi-line): BIBLIO_LIST -> BERTHA (1 line) ---
$d.Paragraphs.Item(19).Range.Text = "2342277 - Bertha María Cuadros Melgar"
